$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.115.25"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "1.791.10"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'224.74"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").Value = "'0.548"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'32.62"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("E9").Value = "  -1.53%  "

$ws.Range("D10").Value = "'0.0706"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "2.049.92"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").Value = "1.790.21"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "'10.81"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").Value = "'0.623"
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").Value = "34.068.13"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").Value = "'68.03"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "'243.26"
$ws.Range("E19").Value = "  -3.26%  "

$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "'10.67"
$ws.Range("E22").Value = "  -3.73%  "

$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("E24").Value = "  -2.71%  "

$ws.Range("D25").Value = "'159.19"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "'7.05"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "'0.0517"
$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("D35").Value = "1.395.58"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  -2.52%  "

$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.35"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "  +2.60%  "

$ws.Range("D41").Value = "'78.92"
$ws.Range("E41").Value = "  -4.55%  "

$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("D43").Value = "'0.914"
$ws.Range("E43").Value = "  -4.34%  "

$ws.Range("E44").Value = "  +18.06%  "

$ws.Range("D45").Value = "'1.07"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").Value = "'108.33"
$ws.Range("E46").Value = "  +2.88%  "

$ws.Range("D47").Value = "'0.0495"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("D49").Value = "1.949.16"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").Value = "'12.12"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.03%  "
